$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.142.64'
$ws.Range("E2").Value = '  -0.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.905.92'
$ws.Range("E3").Value = '  -0.44%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("E4").Value = '  -0.49%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.03'
$ws.Range("E5").Value = '  -0.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4614'
$ws.Range("E7").Value = '  -0.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3890'
$ws.Range("E8").Value = '  -1.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07873'
$ws.Range("E9").Value = '  -0.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9908'
$ws.Range("E10").Value = '  -1.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.98'
$ws.Range("E11").Value = '  -1.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.901.99'
$ws.Range("E12").Value = '  -1.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.756'
$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.041'
$ws.Range("E14").Value = '  -1.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07021'
$ws.Range("E15").Value = '  +0.99%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.14'
$ws.Range("E16").Value = '  -0.49%  '

$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009952'
$ws.Range("E18").Value = '  -1.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.09'
$ws.Range("E19").Value = '  -0.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9998'
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.161.91'
$ws.Range("E21").Value = '  -0.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.322'
$ws.Range("E22").Value = '  -0.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.12'
$ws.Range("E23").Value = '  +0.18%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.092'
$ws.Range("E24").Value = '  +1.46%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.081.42'
$ws.Range("E25").Value = '  -3.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.53'
$ws.Range("E26").Value = '  -0.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.47'
$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.896'
$ws.Range("E28").Value = '  -4.10%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '118.85'
$ws.Range("E29").Value = '  +0.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.879'
$ws.Range("E30").Value = '  -5.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09350'
$ws.Range("E31").Value = '  -0.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8963'
$ws.Range("E32").Value = '  -3.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.233'
$ws.Range("E33").Value = '  -2.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.322'
$ws.Range("E34").Value = '  -2.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.141'
$ws.Range("E35").Value = '  -4.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05786'
$ws.Range("E36").Value = '  -0.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.172'
$ws.Range("E37").Value = '  -2.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02086'
$ws.Range("E38").Value = '  -1.26%  '

$ws.Range("E39").Value = '  -0.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5697'
$ws.Range("E40").Value = '  -0.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.667'
$ws.Range("E41").Value = '  -3.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1809'
$ws.Range("E42").Value = '  +0.29%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.742'
$ws.Range("E43").Value = '  -2.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.94'
$ws.Range("E44").Value = '  -0.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5352'
$ws.Range("E45").Value = '  -1.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.169'
$ws.Range("E46").Value = '  -5.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07002'
$ws.Range("E47").Value = '  -1.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.841'
$ws.Range("E48").Value = '  -2.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.549'
$ws.Range("E49").Value = '  -0.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '113.33'
$ws.Range("E50").Value = '  -0.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2979'
$ws.Range("E51").Value = '  +0.74%  '
